# Daily attendance processing - 2025-10-12 09:18:59
# Applies:
#  - Reorders the "Recorded By" (G column) text for several sessions so
#    that "System" is listed first.
#  - Updates the summary statistics table (K/L columns) to reflect the
#    newly recorded sessions.
#  - Updates the per-group roll-up table (K15:S18) for groups B2A/B2B/B2C.
#  - Marks session 17 (12/10/2025) for groups B2A, B2B and B2C as
#    "Recorded" (was "Pending"), filling in who recorded it and the
#    attendance counts, and re-colors those rows to match the other
#    "Recorded" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a plain-text value into a cell without Excel silently
# re-interpreting it as a number/percentage/date, then restore the
# cell's original formatting (which a raw text write can disturb) by
# pasting the formats from an untouched cell that already has the
# formatting we want to keep.
# ---------------------------------------------------------------------
function Set-TextValue($cellAddr, $text, $formatDonorAddr) {
    $ws.Range($cellAddr).Value = "'" + $text
    if ($formatDonorAddr) {
        $ws.Range($formatDonorAddr).Copy()
        $ws.Range($cellAddr).PasteSpecial(-4122)
    }
}

# ---------------------------------------------------------------------
# "Recorded By" (G column) text reordering - System first
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "System, system, backup@backdoor.com"
$ws.Range("G4").Value = "System, backup@backdoor.com"
$ws.Range("G5").Value = "System, backup@backdoor.com"

$ws.Range("G29").Value = "System, system, backup@backdoor.com"
$ws.Range("G31").Value = "System, backup@backdoor.com"
$ws.Range("G32").Value = "System, backup@backdoor.com"

$ws.Range("G56").Value = "System, system, backup@backdoor.com"
$ws.Range("G58").Value = "System, backup@backdoor.com"
$ws.Range("G59").Value = "System, backup@backdoor.com"

$ws.Range("G84").Value = "System, backup@backdoor.com"
$ws.Range("G85").Value = "System, backup@backdoor.com"

$ws.Range("G110").Value = "System, backup@backdoor.com"
$ws.Range("G111").Value = "System, backup@backdoor.com"

$ws.Range("G136").Value = "System, backup@backdoor.com"
$ws.Range("G137").Value = "System, backup@backdoor.com"

# ---------------------------------------------------------------------
# Overall summary statistics (K3:L10)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 93    # Recorded Sessions (was 90)
$ws.Range("L8").Value = 66    # Pending Sessions (was 69)
Set-TextValue "L9" "58.5%" "L7"     # Coverage % (was 56.6%)
Set-TextValue "L10" "70.4%" "L7"    # Average Attendance % (was 70.8%)

# ---------------------------------------------------------------------
# Per-group roll-up table (K14:S18)
# ---------------------------------------------------------------------
# B2A (row 15)
$ws.Range("O15").Value = 17   # Recorded (was 16)
$ws.Range("Q15").Value = 10   # Pending (was 11)
Set-TextValue "R15" "63.0%" "M15"   # Coverage % (was 59.3%)
Set-TextValue "S15" "69.4%" "M15"   # Avg Attendance % (was 69.6%)

# B2B (row 16)
$ws.Range("O16").Value = 17   # Recorded (was 16)
$ws.Range("Q16").Value = 10   # Pending (was 11)
Set-TextValue "R16" "63.0%" "M16"   # Coverage % (was 59.3%)

# B2C (row 17)
$ws.Range("O17").Value = 17   # Recorded (was 16)
$ws.Range("Q17").Value = 10   # Pending (was 11)
Set-TextValue "R17" "63.0%" "M17"   # Coverage % (was 59.3%)
Set-TextValue "S17" "60.6%" "M17"   # Avg Attendance % (was 61.5%)

# ---------------------------------------------------------------------
# Session 17 (12/10/2025) moves from "Pending" to "Recorded" for
# groups B2A, B2B and B2C. Re-use the formatting of an existing
# "Recorded" row (row 2) so the green highlighting matches.
# ---------------------------------------------------------------------
$ws.Range("A2:I2").Copy()

$ws.Range("A18:I18").PasteSpecial(-4122)
$ws.Range("G18").Value = "dnasr281@gmail.com"
$ws.Range("H18").Value = "35/53"
$ws.Range("I18").Value = "Recorded"

$ws.Range("A45:I45").PasteSpecial(-4122)
$ws.Range("G45").Value = "dnasr281@gmail.com"
$ws.Range("H45").Value = "38/56"
$ws.Range("I45").Value = "Recorded"

$ws.Range("A72:I72").PasteSpecial(-4122)
$ws.Range("G72").Value = "dnasr281@gmail.com"
$ws.Range("H72").Value = "26/55"
$ws.Range("I72").Value = "Recorded"
